$d = $word.ActiveDocument

$replacements = @(
    ,@('6+65=', '6-2=')
    ,@('29+4=', '2+53=')
    ,@('85-35=', '62+13=')
    ,@('85-40=', '41-4=')
    ,@('20+51=', '4+5=')
    ,@('11+75=', '33+42=')
    ,@('2+61=', '62+1=')
    ,@('81+8=', '44+26=')
    ,@('26+63=', '74+17=')
    ,@('92-23=', '77-38=')
    ,@('36+31=', '23+36=')
    ,@('27+70=', '28-15=')
    ,@('38+31=', '61-42=')
    ,@('60+18=', '74-34=')
    ,@('12+9=', '49-18=')
    ,@('74-58=', '82-73=')
    ,@('23+59=', '47-46=')
    ,@('65-25=', '88-57=')
    ,@('83-7=', '42-37=')
    ,@('62-58=', '82-57=')
    ,@('30+5=', '13+30=')
    ,@('50-48=', '6+26=')
    ,@('13+0=', '12+60=')
    ,@('6+39=', '21+27=')
    ,@('41+32=', '24+25=')
    ,@('44-38=', '17-13=')
    ,@('26+57=', '83-12=')
    ,@('43-10=', '47-39=')
    ,@('2+42=', '98-50=')
    ,@('8+6=', '2+55=')
    ,@('42-18=', '51+26=')
    ,@('6+29=', '53+41=')
    ,@('15+53=', '29+49=')
    ,@('2+84=', '52-1=')
    ,@('83-77=', '88+9=')
    ,@('29-17=', '93-70=')
    ,@('29+60=', '77-4=')
    ,@('17-14=', '50+22=')
    ,@('37+29=', '24+47=')
    ,@('97-4=', '16+65=')
    ,@('59-36=', '72-40=')
    ,@('47-35=', '84-13=')
    ,@('94-84=', '60-48=')
    ,@('34-10=', '44-42=')
    ,@('31-3=', '51-24=')
    ,@('10+86=', '91-29=')
    ,@('19+45=', '71+16=')
    ,@('48+44=', '19+48=')
    ,@('75-0=', '96-21=')
    ,@('12+41=', '32+31=')
    ,@('72-8=', '93-50=')
    ,@('1+96=', '51-9=')
    ,@('85-61=', '6-4=')
    ,@('68-37=', '87-11=')
    ,@('82-47=', '76-75=')
    ,@('13+83=', '80+11=')
    ,@('5+87=', '29+64=')
    ,@('73-27=', '97-87=')
    ,@('77-58=', '33-9=')
    ,@('72-70=', '66+20=')
    ,@('17+27=', '44+53=')
    ,@('70-29=', '73-22=')
    ,@('65-54=', '8+71=')
    ,@('94-40=', '52-15=')
    ,@('94-15=', '58-2=')
    ,@('94-93=', '6+60=')
    ,@('14+59=', '48-14=')
    ,@('13+21=', '97-47=')
    ,@('83-9=', '84-81=')
    ,@('38-14=', '52+25=')
    ,@('5+8=', '90-28=')
    ,@('60+12=', '99-66=')
    ,@('2+30=', '84-37=')
    ,@('22+17=', '24+49=')
    ,@('59+25=', '49-38=')
    ,@('81-56=', '89+7=')
    ,@('67-61=', '73-12=')
    ,@('51-22=', '44+53=')
    ,@('54+5=', '95-32=')
    ,@('58-32=', '94-44=')
    ,@('40-24=', '0+89=')
    ,@('78-30=', '55-3=')
    ,@('76-19=', '53+24=')
    ,@('22-20=', '19+73=')
    ,@('37-11=', '82-40=')
    ,@('25+17=', '2+57=')
    ,@('37+45=', '70+25=')
    ,@('70-24=', '18-13=')
    ,@('66-18=', '88-82=')
    ,@('64-53=', '79-5=')
    ,@('42+37=', '66-44=')
    ,@('25+24=', '44+9=')
    ,@('14+60=', '67-63=')
    ,@('87-76=', '32+42=')
    ,@('26-0=', '35+6=')
    ,@('70+13=', '70-46=')
    ,@('21+25=', '19+5=')
    ,@('79-66=', '34+15=')
    ,@('46+28=', '95-46=')
    ,@('89-6=', '93-22=')
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Write-Output "Done applying replacements"
